$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-coerced to numbers by Excel (the source data keeps these as text).
# (Multi-area "A1,A2,A3" Range strings only affect the first area here,
# so set NumberFormat on each cell individually.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.521.70"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").Value = "2.061.61"
$ws.Range("E3").Value = "  +3.99%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "253.81"
$ws.Range("E5").Value = "  +3.00%  "
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("D7").Value = "67.94"
$ws.Range("E7").Value = "  +16.50%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +6.93%  "
$ws.Range("D10").Value = "59.97"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").Value = "0.0772"
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "0.937"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "23.68"
$ws.Range("E14").Value = "  +28.61%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "14.99"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "2.359.42"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "5.62"
$ws.Range("E17").Value = "  +6.04%  "
$ws.Range("D18").Value = "2.061.99"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("D19").Value = "37.473.04"
$ws.Range("E19").Value = "  +5.49%  "
$ws.Range("D20").Value = "73.78"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("D21").Value = "0.0₃0881"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").Value = "5.50"
$ws.Range("E22").Value = "  +5.39%  "
$ws.Range("D23").Value = "240.24"
$ws.Range("D24").Value = "2.74"
$ws.Range("E24").Value = "  +5.78%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +10.08%  "
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  +9.69%  "
$ws.Range("D28").Value = "162.46"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "20.13"
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +38.39%  "
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("E32").Value = "  +7.81%  "
$ws.Range("E33").Value = "  +9.10%  "
$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  +8.58%  "
$ws.Range("D35").Value = "0.0631"
$ws.Range("E35").Value = "  +6.17%  "
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "6.25"
$ws.Range("E37").Value = "  +15.75%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +38.92%  "
$ws.Range("E41").Value = "  +14.70%  "
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").Value = "3.05"
$ws.Range("E43").Value = "  +5.94%  "
$ws.Range("D44").Value = "17.64"
$ws.Range("E44").Value = "  +9.14%  "
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").Value = "97.77"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("D48").Value = "8.00"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").Value = "1.420.38"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "3.76"
$ws.Range("E51").Value = "  +10.18%  "
